$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update Price (D) / Volume(1h) (E) columns,
# and rotate the WrappedEther/TRON/Solana rows (11-13) to their new rank order.
#
# Some new Price values read like plain numbers (e.g. "314.25", "1.003") even
# though the column stores them as text (dotted price strings, leading zeros,
# etc.). Those are written with a leading apostrophe so Excel keeps them as
# literal text instead of silently converting them to numbers, then the style
# is reset to "Normal" so the cell ends up plain text/General - exactly like
# the original cell - without leaving a stray "number stored as text" marker.

$ws.Range('D2').Value = '27.334.68'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '1.856.83'
$ws.Range('E3').Value = '  +1.85%  '
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').Value = '''314.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '''0.4635'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').Value = '''0.3717'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').Value = '''0.07343'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.59%  '
$ws.Range('D10').Value = '''0.8825'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.73%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '''0.07897'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.68%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '''19.89'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.960.58'
$ws.Range('E13').Value = '  +4.50%  '
$ws.Range('D14').Value = '''5.395'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '''6.569'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').Value = '''92.03'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').Value = '''0.000008878'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.63%  '
$ws.Range('D19').Value = '''1.003'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').Value = '''14.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('D21').Value = '27.370.10'
$ws.Range('E21').Value = '  +2.26%  '
$ws.Range('D22').Value = '''5.132'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '2.136.87'
$ws.Range('E24').Value = '  -0.99%  '
$ws.Range('D25').Value = '''152.81'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').Value = '''1.885'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.54%  '
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('D28').Value = '''2.082'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('D29').Value = '''5.124'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('D30').Value = '''116.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.68%  '
$ws.Range('D31').Value = '''0.08889'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').Value = '''0.7571'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.17%  '
$ws.Range('D33').Value = '''3.021'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('D34').Value = '''1.164'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.10%  '
$ws.Range('D35').Value = '''4.488'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('D36').Value = '''2.609'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.98%  '
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = '''0.01955'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('D39').Value = '''2.976'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('D40').Value = '''0.05228'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('D41').Value = '''7.098'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').Value = '''0.5162'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D44').Value = '''8.334'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('D45').Value = '''0.4844'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('D46').Value = '''10.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').Value = '''103.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').Value = '''1.656'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '''65.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.76%  '
